# Reorders the "Requisitos" bullet list so that the "LOB1019 - Física II"
# requirement line becomes the first entry instead of the last one, while
# leaving the other two requirement lines (and their own runs/line breaks)
# untouched.

$d = $word.ActiveDocument

# Locate the "Requisitos" heading paragraph, then grab the very next
# paragraph - that's the ListBullet paragraph holding the requirement runs.
$targetParagraph = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Requisitos") {
        $targetParagraph = $para.Next()
        break
    }
}
if ($targetParagraph -eq $null) {
    throw "Could not find the 'Requisitos' bullet-list paragraph."
}

$newLine = [char]11
$movedText = "LOB1019 -  Física II  (Requisito fraco)"

$paraRange = $targetParagraph.Range

# Insert the moved requirement (text + manual line break) as a brand new
# run right at the start of the paragraph.
$insertionPoint = $d.Range($paraRange.Start, $paraRange.Start)
$insertionPoint.InsertBefore($movedText + $newLine)

# Re-fetch the paragraph range (it grew) and find the *second* occurrence
# of the moved requirement's text+break - i.e. the original run that is
# still sitting at the end of the paragraph - and delete it.
$refreshedParagraph = $targetParagraph.Range
$searchRange = $d.Range($refreshedParagraph.Start + 1, $refreshedParagraph.End)
$found = $searchRange.Find.Execute($movedText + $newLine, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $searchRange.Text = ""
}
